$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(32, 8).Value = 432
$ws.Cells.Item(32, 10).Value = 432
$ws.Cells.Item(32, 12).Value = 432
$ws.Cells.Item(32, 14).Value = -1084
$ws.Cells.Item(137, 8).Value = 2495.8276
$ws.Cells.Item(137, 9).Value = 1222.6154
$ws.Cells.Item(137, 10).Value = 3530.3125
$ws.Cells.Item(137, 11).Value = 3667.8462
$ws.Cells.Item(137, 12).Value = 10590.9375
$ws.Cells.Item(137, 13).Value = -1117.8462
$ws.Cells.Item(137, 14).Value = -15690.9375
$ws.Cells.Item(138, 8).Value = 2993.29
$ws.Cells.Item(138, 9).Value = 1399.375
$ws.Cells.Item(138, 10).Value = 3743.3677
$ws.Cells.Item(138, 11).Value = 4198.125
$ws.Cells.Item(138, 12).Value = 11230.1031
$ws.Cells.Item(138, 13).Value = 941.875
$ws.Cells.Item(138, 14).Value = -21510.1031

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 2625.1875
$ws.Cells.Item(61, 9).Value = 2466.8667
$ws.Cells.Item(61, 10).Value = 5000
$ws.Cells.Item(61, 11).Value = 2466.8667
$ws.Cells.Item(61, 12).Value = 5000
$ws.Cells.Item(61, 13).Value = -2254.8667
$ws.Cells.Item(61, 14).Value = -5424
$ws.Cells.Item(80, 8).Value = 29403.334
$ws.Cells.Item(80, 9).Value = 9100
$ws.Cells.Item(80, 10).Value = 39555
$ws.Cells.Item(80, 11).Value = 9100
$ws.Cells.Item(80, 12).Value = 39555
$ws.Cells.Item(80, 13).Value = -8102
$ws.Cells.Item(80, 14).Value = -41551
$ws.Cells.Item(83, 8).Value = 29403.334
$ws.Cells.Item(83, 9).Value = 9100
$ws.Cells.Item(83, 10).Value = 39555
$ws.Cells.Item(83, 11).Value = 27300
$ws.Cells.Item(83, 12).Value = 118665
$ws.Cells.Item(83, 13).Value = -22308
$ws.Cells.Item(83, 14).Value = -128649
$ws.Cells.Item(132, 8).Value = 1975.425
$ws.Cells.Item(132, 9).Value = 1546.2858
$ws.Cells.Item(132, 10).Value = 4979.4
$ws.Cells.Item(132, 11).Value = 4638.857400000001
$ws.Cells.Item(132, 12).Value = 14938.2
$ws.Cells.Item(132, 13).Value = -2108.857400000001
$ws.Cells.Item(132, 14).Value = -19998.2
$ws.Cells.Item(136, 8).Value = 2625.1875
$ws.Cells.Item(136, 9).Value = 2466.8667
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 7400.6001
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = -4850.6001
$ws.Cells.Item(136, 14).Value = -20100
$ws.Cells.Item(139, 8).Value = 145099.8
$ws.Cells.Item(139, 10).Value = 145099.8
$ws.Cells.Item(139, 12).Value = 145099.8
$ws.Cells.Item(139, 14).Value = -155379.8

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(82, 8).Value = 16247.211
$ws.Cells.Item(82, 9).Value = 3717.5833
$ws.Cells.Item(82, 10).Value = 37726.57
$ws.Cells.Item(82, 11).Value = 3717.5833
$ws.Cells.Item(82, 12).Value = 37726.57
$ws.Cells.Item(82, 13).Value = -3334.5833
$ws.Cells.Item(82, 14).Value = -38492.57
$ws.Cells.Item(85, 8).Value = 16247.211
$ws.Cells.Item(85, 9).Value = 3717.5833
$ws.Cells.Item(85, 10).Value = 37726.57
$ws.Cells.Item(85, 11).Value = 3717.5833
$ws.Cells.Item(85, 12).Value = 37726.57
$ws.Cells.Item(85, 13).Value = -2391.5833
$ws.Cells.Item(85, 14).Value = -40378.57

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 8).Value = 20800.4
$ws.Cells.Item(4, 9).Value = 1000
$ws.Cells.Item(4, 10).Value = 100002
$ws.Cells.Item(4, 11).Value = 1000
$ws.Cells.Item(4, 12).Value = 100002
$ws.Cells.Item(4, 13).Value = -888
$ws.Cells.Item(4, 14).Value = -100226
$ws.Cells.Item(31, 8).Value = 3854.3076
$ws.Cells.Item(31, 9).Value = 2009.92
$ws.Cells.Item(31, 10).Value = 7147.857
$ws.Cells.Item(31, 11).Value = 2009.92
$ws.Cells.Item(31, 12).Value = 7147.857
$ws.Cells.Item(31, 13).Value = -1714.92
$ws.Cells.Item(31, 14).Value = -7737.857
$ws.Cells.Item(34, 8).Value = 3854.3076
$ws.Cells.Item(34, 9).Value = 2009.92
$ws.Cells.Item(34, 10).Value = 7147.857
$ws.Cells.Item(34, 11).Value = 2009.92
$ws.Cells.Item(34, 12).Value = 7147.857
$ws.Cells.Item(34, 13).Value = -1807.92
$ws.Cells.Item(34, 14).Value = -7551.857
$ws.Cells.Item(41, 8).Value = 17917.6
$ws.Cells.Item(41, 10).Value = 24843.334
$ws.Cells.Item(41, 12).Value = 24843.334
$ws.Cells.Item(41, 14).Value = -25699.334
$ws.Cells.Item(50, 8).Value = 9222.666999999999
$ws.Cells.Item(50, 10).Value = 9222.666999999999
$ws.Cells.Item(50, 12).Value = 9222.666999999999
$ws.Cells.Item(50, 14).Value = -10472.667
$ws.Cells.Item(51, 8).Value = 9249.714
$ws.Cells.Item(51, 10).Value = 9249.714
$ws.Cells.Item(51, 12).Value = 9249.714
$ws.Cells.Item(51, 14).Value = -10721.714
$ws.Cells.Item(60, 8).Value = 23749.215
$ws.Cells.Item(60, 10).Value = 23749.215
$ws.Cells.Item(60, 12).Value = 23749.215
$ws.Cells.Item(60, 14).Value = -24771.215
$ws.Cells.Item(61, 8).Value = 9249.714
$ws.Cells.Item(61, 10).Value = 9249.714
$ws.Cells.Item(61, 12).Value = 9249.714
$ws.Cells.Item(61, 14).Value = -9945.714
$ws.Cells.Item(68, 8).Value = 17449.5
$ws.Cells.Item(68, 10).Value = 17449.5
$ws.Cells.Item(68, 12).Value = 17449.5
$ws.Cells.Item(68, 14).Value = -18947.5
$ws.Cells.Item(71, 8).Value = 17449.5
$ws.Cells.Item(71, 10).Value = 17449.5
$ws.Cells.Item(71, 12).Value = 52348.5
$ws.Cells.Item(71, 14).Value = -59836.5
$ws.Cells.Item(109, 8).Value = 18185
$ws.Cells.Item(109, 10).Value = 18185
$ws.Cells.Item(109, 12).Value = 18185
$ws.Cells.Item(109, 14).Value = -20265
$ws.Cells.Item(122, 8).Value = 1875.6
$ws.Cells.Item(122, 9).Value = 1417.3334
$ws.Cells.Item(122, 11).Value = 4252.0002
$ws.Cells.Item(122, 13).Value = -1802.0002

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(121, 8).Value = 2669226.5
$ws.Cells.Item(121, 10).Value = 6671666.5
$ws.Cells.Item(121, 12).Value = 20014999.5
$ws.Cells.Item(121, 14).Value = -20017619.5

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(20, 8).Value = 20581
$ws.Cells.Item(20, 9).Value = 23952.5
$ws.Cells.Item(20, 10).Value = 18333.334
$ws.Cells.Item(20, 11).Value = 23952.5
$ws.Cells.Item(20, 12).Value = 18333.334
$ws.Cells.Item(20, 13).Value = -23726.5
$ws.Cells.Item(20, 14).Value = -18785.334
$ws.Cells.Item(46, 8).Value = 2266.8
$ws.Cells.Item(46, 9).Value = 1360
$ws.Cells.Item(46, 10).Value = 2720.2
$ws.Cells.Item(46, 11).Value = 1360
$ws.Cells.Item(46, 12).Value = 2720.2
$ws.Cells.Item(46, 13).Value = -1172
$ws.Cells.Item(46, 14).Value = -3096.2
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(132, 8).Value = 4184.3887
$ws.Cells.Item(132, 9).Value = 4729.5454
$ws.Cells.Item(132, 10).Value = 3327.7144
$ws.Cells.Item(132, 11).Value = 14188.6362
$ws.Cells.Item(132, 12).Value = 9983.143199999999
$ws.Cells.Item(132, 13).Value = -11658.6362
$ws.Cells.Item(132, 14).Value = -15043.1432

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(129, 8).Value = 42500
$ws.Cells.Item(129, 10).Value = 42500
$ws.Cells.Item(129, 12).Value = 42500
$ws.Cells.Item(129, 14).Value = -52500

$ws = $wb.Worksheets.Item(7)
$ws.Range("N109").ClearContents()
$ws.Range("N115").ClearContents()
